$d = $word.ActiveDocument
$d.TrackRevisions = $false

# 1. Update deadline date: "November 20, 2022" -> "November 19, 2023"
$r = $d.Content
$r.Find.Execute("November 20, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "November 19, 2023", 2) | Out-Null

# 2. "G(0,t)" -> "Gt(0)"  (keep italics formatting of the surrounding run)
$r = $d.Content
$r.Find.Execute("G(0,t)", $true, $false, $false, $false, $false, $true, 1, $false, "Gt(0)", 2) | Out-Null

# 3. Subscript "t" in "Kt" -> "T"  (only change the subscript character, preserve run formatting)
$r = $d.Content
$found = $r.Find.Execute("Kt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $sub = $d.Range($r.Start + 1, $r.End)
    $sub.Text = "T"
}

# 4. "model the complete year." -> "model the entire year."
$r = $d.Content
$r.Find.Execute("model the complete year.", $true, $false, $false, $false, $false, $true, 1, $false, "model the entire year.", 2) | Out-Null

# 5. Trim the "Note:" sentence
$r = $d.Content
$r.Find.Execute("Note: Solar position and irradiance equations can be found at the lecture slides, as well as in Appendix A of the paper provided in the section Further Readings of the lecture VII. Solar radiation. ", $true, $false, $false, $false, $false, $true, 1, $false, "Note: Solar position and irradiance equations can be found at the lecture slides.", 2) | Out-Null

# 6. "Let's assume" -> "Let us assume"
$r = $d.Content
$r.Find.Execute("Let" + [char]8217 + "s assume", $true, $false, $false, $false, $false, $true, 1, $false, "Let us assume", 2) | Out-Null

# 7. "power decrease due to ambient temperature" -> "power reduction due to ambient temperature"
$r = $d.Content
$r.Find.Execute("power decrease due to ambient temperature", $true, $false, $false, $false, $false, $true, 1, $false, "power reduction due to ambient temperature", 2) | Out-Null
